$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REMISIONES OCTUBRE  2021     ")

# Row 14 - complete the payment date/amount (credit from 22-Oct closes)
$ws.Range("F14").Value = 44492
$ws.Range("G14").Value = 5290

# Row 15 - new credit entry 24-Oct, DAVID
$ws.Range("A15").Value = 44493
$ws.Range("D15").Value = "DAVID"
$ws.Range("E15").Value = 3208
$ws.Range("F15").Value = 44493
$ws.Range("G15").Value = 3208

# Row 16 - new credit entry 24-Oct, GUSTAVO
$ws.Range("A16").Value = 44493
$ws.Range("D16").Value = "GUSTAVO"
$ws.Range("E16").Value = 3500
$ws.Range("F16").Value = 44497
$ws.Range("G16").Value = 3500

# Row 17 - new credit entry 25-Oct, HERRADURA DAVID
$ws.Range("A17").Value = 44494
$ws.Range("D17").Value = "HERRADURA DAVID"
$ws.Range("E17").Value = 5511
$ws.Range("F17").Value = 44495
$ws.Range("G17").Value = 5511

# Row 18 - new credit entry 28-Oct, GUSTAVO (still unpaid, F/G left blank)
$ws.Range("A18").Value = 44497
$ws.Range("D18").Value = "GUSTAVO"
$ws.Range("E18").Value = 3511

# Update the selected cell to reflect where the user left off
[void]$ws.Range("G17").Select()
